# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the Titan_Profits.xlsx diff (refreshed market-price derived columns).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 649.0909
$ws.Range("I6").Value = 450
$ws.Range("J6").Value = 997.5
$ws.Range("K6").Value = 1350
$ws.Range("L6").Value = 2992.5
$ws.Range("M6").Value = -1238
$ws.Range("N6").Value = -3216.5

# Row 9
$ws.Range("H9").Value = 85.166664
$ws.Range("I9").Value = 92.59999999999999
$ws.Range("J9").Value = 48
$ws.Range("K9").Value = 92.59999999999999
$ws.Range("L9").Value = 48
$ws.Range("M9").Value = 76.40000000000001
$ws.Range("N9").Value = -386

# Row 12
$ws.Range("H12").Value = 51.75
$ws.Range("I12").Value = 68.333336
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 68.333336
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 101.666664

# Row 28
$ws.Range("H28").Value = 742063.0600000001
$ws.Range("I28").Value = 1235440.6
$ws.Range("J28").Value = 1996.6666
$ws.Range("K28").Value = 1235440.6
$ws.Range("L28").Value = 1996.6666
$ws.Range("M28").Value = -1234955.6

# Row 33
$ws.Range("H33").Value = 318.75
$ws.Range("I33").Value = 550
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 550
$ws.Range("L33").Value = 180
$ws.Range("M33").Value = -321
$ws.Range("N33").Value = -638

# Row 38
$ws.Range("H38").Value = 18.6
$ws.Range("I38").Value = 18.6
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 55.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 316.2

# Row 70
$ws.Range("H70").Value = 2497.7144
$ws.Range("I70").Value = 1798.75
$ws.Range("J70").Value = 2662.1765
$ws.Range("K70").Value = 5396.25
$ws.Range("L70").Value = 7986.529500000001
$ws.Range("M70").Value = -5126.25
$ws.Range("N70").Value = -8526.529500000001

# Row 73
$ws.Range("H73").Value = 2497.7144
$ws.Range("I73").Value = 1798.75
$ws.Range("J73").Value = 2662.1765
$ws.Range("K73").Value = 5396.25
$ws.Range("L73").Value = 7986.529500000001
$ws.Range("M73").Value = -4460.25
$ws.Range("N73").Value = -9858.529500000001

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

# Row 125
$ws.Range("H125").Value = 7938159.5
$ws.Range("I125").Value = 810.5
$ws.Range("J125").Value = 13891171
$ws.Range("K125").Value = 7294.5
$ws.Range("L125").Value = 125020539
$ws.Range("M125").Value = -4834.5
$ws.Range("N125").Value = -125025459

# Row 138
$ws.Range("H138").Value = 6423986
$ws.Range("I138").Value = 3108826
$ws.Range("J138").Value = 7579269
$ws.Range("K138").Value = 9326478
$ws.Range("L138").Value = 22737807
$ws.Range("M138").Value = -9321338
$ws.Range("N138").Value = -22748087

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 127541.625
$ws.Range("I2").Value = 145590.42
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 145590.42
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -145477.42
$ws.Range("N2").Value = -1426

# Row 74
$ws.Range("H74").Value = 5229.364
$ws.Range("I74").Value = 1097.8572
$ws.Range("J74").Value = 12459.5
$ws.Range("K74").Value = 1097.8572
$ws.Range("L74").Value = 12459.5
$ws.Range("M74").Value = -223.8571999999999
$ws.Range("N74").Value = -14207.5

# Row 77
$ws.Range("H77").Value = 5229.364
$ws.Range("I77").Value = 1097.8572
$ws.Range("J77").Value = 12459.5
$ws.Range("K77").Value = 5489.286
$ws.Range("L77").Value = 62297.5
$ws.Range("M77").Value = -1121.286
$ws.Range("N77").Value = -71033.5

# Row 116
$ws.Range("H116").Value = 127541.625
$ws.Range("I116").Value = 145590.42
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 145590.42
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = -143296.42
$ws.Range("N116").Value = -5788

# Row 132
$ws.Range("H132").Value = 3721.6
$ws.Range("I132").Value = 3075.9092
$ws.Range("J132").Value = 5497.25
$ws.Range("K132").Value = 9227.7276
$ws.Range("L132").Value = 16491.75
$ws.Range("M132").Value = -6697.7276
$ws.Range("N132").Value = -21551.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 127541.625
$ws.Range("I3").Value = 145590.42
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 145590.42
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = -145476.42
$ws.Range("N3").Value = -1428

# Row 20
$ws.Range("H20").Value = 3799.8
$ws.Range("I20").Value = 3999.6667
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 3999.6667
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -3752.6667
$ws.Range("N20").Value = -3994

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 105
$ws.Range("H105").Value = 3392.1
$ws.Range("I105").Value = 3392.2856
$ws.Range("J105").Value = 3391.6667
$ws.Range("K105").Value = 3392.2856
$ws.Range("L105").Value = 3391.6667
$ws.Range("M105").Value = -1645.2856
$ws.Range("N105").Value = -6885.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 895.06665
$ws.Range("I105").Value = 754.43475
$ws.Range("J105").Value = 1357.1428
$ws.Range("K105").Value = 754.43475
$ws.Range("L105").Value = 1357.1428
$ws.Range("M105").Value = 992.56525

# Row 134
$ws.Range("H134").Value = 40543468
$ws.Range("I134").Value = 66668296
$ws.Range("J134").Value = 22731088
$ws.Range("K134").Value = 200004888
$ws.Range("L134").Value = 68193264
$ws.Range("M134").Value = -200002353
$ws.Range("N134").Value = -68198334

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 1142.7142
$ws.Range("I80").Value = 999
$ws.Range("J80").Value = 1166.6666
$ws.Range("K80").Value = 2997
$ws.Range("L80").Value = 3499.9998
$ws.Range("M80").Value = -2061
$ws.Range("N80").Value = -5371.9998

# Row 83
$ws.Range("H83").Value = 1142.7142
$ws.Range("I83").Value = 999
$ws.Range("J83").Value = 1166.6666
$ws.Range("K83").Value = 8991
$ws.Range("L83").Value = 10499.9994
$ws.Range("M83").Value = -4311
$ws.Range("N83").Value = -19859.9994

# Row 132
$ws.Range("H132").Value = 1526.6666
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 1585.7142
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 14271.4278
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -19331.4278

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 58826364
$ws.Range("I80").Value = 2709.0908
$ws.Range("J80").Value = 166669730
$ws.Range("K80").Value = 2709.0908
$ws.Range("L80").Value = 166669730
$ws.Range("M80").Value = -1711.0908
$ws.Range("N80").Value = -166671726

# Row 83
$ws.Range("H83").Value = 58826364
$ws.Range("I83").Value = 2709.0908
$ws.Range("J83").Value = 166669730
$ws.Range("K83").Value = 13545.454
$ws.Range("L83").Value = 833348650
$ws.Range("M83").Value = -8553.454
$ws.Range("N83").Value = -833358634

$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Range("H24").Value = 4500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 4500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 4500
$ws.Range("N24").Value = -5186

# Row 25
$ws.Range("H25").Value = 50000
$ws.Range("I25").Value = 50000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 50000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -49770
$ws.Range("N25").ClearContents()

# Row 36
$ws.Range("H36").Value = 26890.572
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 26890.572
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 26890.572
$ws.Range("N36").Value = -28014.572

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 9374.75
$ws.Range("I20").Value = 7499
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 7499
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -7259
$ws.Range("N20").Value = -10480

# Row 81
$ws.Range("H81").Value = 4121.6523
$ws.Range("I81").Value = 825
$ws.Range("J81").Value = 4815.684
$ws.Range("K81").Value = 1650
$ws.Range("L81").Value = 9631.368
$ws.Range("M81").Value = -589
$ws.Range("N81").Value = -11753.368

# Row 84
$ws.Range("H84").Value = 4121.6523
$ws.Range("I84").Value = 825
$ws.Range("J84").Value = 4815.684
$ws.Range("K84").Value = 8250
$ws.Range("L84").Value = 48156.84
$ws.Range("M84").Value = -2946
$ws.Range("N84").Value = -58764.84

# Row 122
$ws.Range("H122").Value = 3100
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9300
$ws.Range("N122").Value = -14200
